# Weighting & Scaling update & heatmap
$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("ecological_params")
$wsScaling = $wb.Worksheets.Item("Scaling")

# --- Scaling sheet: add new header columns E/F ---
$wsScaling.Range("E1").Value = "Optimal"
$wsScaling.Range("F1").Value = "Threshold"
$wsScaling.Range("C1").Copy()
$wsScaling.Range("E1:F1").PasteSpecial(-4122) # xlPasteFormats - match header formatting

# --- Scaling sheet: update Min column (B) to 0 ---
$wsScaling.Range("B2").Value = 0
$wsScaling.Range("B3").Value = 0
$wsScaling.Range("B4").Value = 0

# --- Scaling sheet: update Max column (C) to reference MAX of ecological_params row ---
$wsScaling.Range("C2").Formula = "=MAX(ecological_params!B2:D2)"
$wsScaling.Range("C3").Formula = "=MAX(ecological_params!B3:D3)"
$wsScaling.Range("C4").Formula = "=MAX(ecological_params!B4:D4)"

# --- Selections to reflect final cursor/viewport state ---
$wsScaling.Select()
$wsScaling.Range("C2:C4").Select()

$wsParams.Select()
$wsParams.Range("C13").Select()

$wb.Save()
